$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.098.81'
$ws.Range("E2").Value = '  +5.31%  '

$ws.Range("D3").Value = '3.371.17'
$ws.Range("E3").Value = '  +5.87%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''573.59'
$ws.Range("E5").Value = '  +7.43%  '

$ws.Range("D6").Value = '''152.85'
$ws.Range("E6").Value = '  +5.47%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.382.56'
$ws.Range("E8").Value = '  +5.96%  '

$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").Value = '''7.42'
$ws.Range("E10").Value = '  +1.56%  '

$ws.Range("E11").Value = '  +6.48%  '

$ws.Range("E12").Value = '  +1.99%  '

$ws.Range("D13").Value = '3.948.42'
$ws.Range("E13").Value = '  +5.76%  '

$ws.Range("E14").Value = '  +0.32%  '

$ws.Range("D15").Value = '''27.00'
$ws.Range("E15").Value = '  +4.63%  '

$ws.Range("D16").Value = '''0.0000181'
$ws.Range("E16").Value = '  +4.99%  '

$ws.Range("D17").Value = '63.098.71'
$ws.Range("E17").Value = '  +5.23%  '

$ws.Range("D18").Value = '3.367.03'
$ws.Range("E18").Value = '  +5.35%  '

$ws.Range("D19").Value = '''6.32'
$ws.Range("E19").Value = '  +1.39%  '

$ws.Range("D20").Value = '''13.92'
$ws.Range("E20").Value = '  +5.47%  '

$ws.Range("D21").Value = '''8.42'
$ws.Range("E21").Value = '  +3.01%  '

$ws.Range("D22").Value = '''385.00'
$ws.Range("E22").Value = '  +4.91%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("E24").Value = '  +2.77%  '

$ws.Range("D25").Value = '''70.48'
$ws.Range("E25").Value = '  +1.45%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = '''0.179'
$ws.Range("E26").Value = '  +6.78%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '''9.22'
$ws.Range("E27").Value = '  +7.13%  '

$ws.Range("D28").Value = '0.0₃0969'
$ws.Range("E28").Value = '  +12.01%  '

$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("E30").Value = '  +6.92%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''23.06'
$ws.Range("E31").Value = '  +3.55%  '

$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '''1.32'
$ws.Range("E32").Value = '  +11.43%  '

$ws.Range("D33").Value = '''5.56'
$ws.Range("E33").Value = '  +5.79%  '

$ws.Range("D34").Value = '''6.30'
$ws.Range("E34").Value = '  +4.48%  '

$ws.Range("D35").Value = '''6.72'
$ws.Range("E35").Value = '  +2.61%  '

$ws.Range("E36").Value = '  +9.65%  '

$ws.Range("D37").Value = '''158.19'
$ws.Range("E37").Value = '  +1.74%  '

$ws.Range("E38").Value = '  +12.47%  '

$ws.Range("D39").Value = '''27.40'
$ws.Range("E39").Value = '  +5.49%  '

$ws.Range("D40").Value = '2.889.96'
$ws.Range("E40").Value = '  +2.81%  '

$ws.Range("D41").Value = '''0.0329'
$ws.Range("E41").Value = '  +11.09%  '

$ws.Range("E42").Value = '  +5.79%  '

$ws.Range("D43").Value = '''40.79'
$ws.Range("E43").Value = '  +3.06%  '

$ws.Range("D44").Value = '''0.746'
$ws.Range("E44").Value = '  +4.32%  '

$ws.Range("E45").Value = '  +1.20%  '

$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '''1.04'
$ws.Range("E46").Value = '  +6.10%  '

$ws.Range("B47").Value = 'RenzoRestakedETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D47").Value = '3.412.85'
$ws.Range("E47").Value = '  +5.76%  '

$ws.Range("D48").Value = '''21.92'
$ws.Range("E48").Value = '  +6.91%  '

$ws.Range("D49").Value = '''299.91'
$ws.Range("E49").Value = '  +14.09%  '

$ws.Range("E50").Value = '  -1.65%  '

$ws.Range("D51").Value = '''6.30'
$ws.Range("E51").Value = '  +2.56%  '

